# Updated cryptos list values (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.361.40"
$ws.Range("E2").Value = "  -7.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.678.21"
$ws.Range("E3").Value = "  -5.88%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.48"
$ws.Range("E5").Value = "  -4.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5094"
$ws.Range("E6").Value = "  -12.85%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2657"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "22.07"
$ws.Range("E9").Value = "  -4.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06326"
$ws.Range("E10").Value = "  -4.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07358"
$ws.Range("E11").Value = "  -2.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.675.42"
$ws.Range("E12").Value = "  -6.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.542"
$ws.Range("E13").Value = "  -4.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5759"
$ws.Range("E14").Value = "  -4.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.908.03"
$ws.Range("E15").Value = "  -5.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008564"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.79"
$ws.Range("E17").Value = "  -13.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.414.08"
$ws.Range("E18").Value = "  -6.96%  "
$ws.Range("E19").Value = "  -6.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("E21").Value = "  -4.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "186.04"
$ws.Range("E22").Value = "  -9.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.230"
$ws.Range("E23").Value = "  -7.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.78"
$ws.Range("E25").Value = "  -5.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.533"
$ws.Range("E26").Value = "  -6.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1171"
$ws.Range("E27").Value = "  -5.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.78"
$ws.Range("E28").Value = "  -2.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.331"
$ws.Range("E29").Value = "  -5.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05850"
$ws.Range("E30").Value = "  -3.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.333"
$ws.Range("E31").Value = "  -5.65%  "
$ws.Range("E32").Value = "  -6.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.503"
$ws.Range("E33").Value = "  -6.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.655"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.007"
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5950"
$ws.Range("E36").Value = "  -6.41%  "
$ws.Range("E37").Value = "  -5.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.672"
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.096.79"
$ws.Range("E40").Value = "  -3.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.894"
$ws.Range("E41").Value = "  -5.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8600"
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.81"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.834.56"
$ws.Range("E45").Value = "  -5.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000113"
$ws.Range("E46").Value = "  +4.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.26"
$ws.Range("E47").Value = "  -5.74%  "
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.051"
$ws.Range("E49").Value = "  -3.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05203"
$ws.Range("E51").Value = "  -3.90%  "
